$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1 (13:29 -> 13:59)
$ws.Range("A1").Value = "Datos actualizados a 28 de Marzo de 2020 a las 13:59"

# Refresh country case numbers and re-sort by Casos totales (column B) descending.
# Only cells whose resulting value actually changes are written below.

# Row 16: Austria -> Austria
$ws.Range("B16").Value = 7964
$ws.Range("C16").Value = 267
$ws.Range("E16").Value = 7671

# Row 19: Portugal -> Portugal
$ws.Range("E19").Value = 4125
$ws.Range("G19").Value = 24
$ws.Range("H19").Value = 100

# Row 35: Pakistan -> Pakistan
$ws.Range("B35").Value = 1415
$ws.Range("C35").Value = 42
$ws.Range("D35").Value = 29
$ws.Range("E35").Value = 1374
$ws.Range("G35").Value = 1
$ws.Range("H35").Value = 12

# Row 38: Sudafrica -> Arabia Saudita
$ws.Range("A38").Value = "Arabia Saudita"
$ws.Range("B38").Value = 1203
$ws.Range("C38").Value = 99
$ws.Range("D38").Value = 37
$ws.Range("E38").Value = 1162
$ws.Range("F38").Value = 6
$ws.Range("G38").Value = 1
$ws.Range("H38").Value = 4

# Row 39: Finlandia -> Sudafrica
$ws.Range("A39").Value = "Sudafrica"
$ws.Range("B39").Value = 1170
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 31
$ws.Range("E39").Value = 1138
$ws.Range("F39").Value = 7
$ws.Range("H39").Value = 1

# Row 40: Indonesia -> Finlandia
$ws.Range("A40").Value = "Finlandia"
$ws.Range("B40").Value = 1167
$ws.Range("C40").Value = 126
$ws.Range("D40").Value = 10
$ws.Range("E40").Value = 1149
$ws.Range("F40").Value = 32
$ws.Range("G40").Value = 1
$ws.Range("H40").Value = 8

# Row 41: Arabia Saudita -> Indonesia
$ws.Range("A41").Value = "Indonesia"
$ws.Range("B41").Value = 1155
$ws.Range("C41").Value = 109
$ws.Range("D41").Value = 59
$ws.Range("E41").Value = 994
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 15
$ws.Range("H41").Value = 102

# Row 57: Colombia -> Hong Kong
$ws.Range("A57").Value = "Hong Kong"
$ws.Range("B57").Value = 560
$ws.Range("C57").Value = 41
$ws.Range("D57").Value = 112
$ws.Range("E57").Value = 444
$ws.Range("F57").Value = 5
$ws.Range("H57").Value = 4

# Row 58: Egipto -> Colombia
$ws.Range("A58").Value = "Colombia"
$ws.Range("B58").Value = 539
$ws.Range("D58").Value = 10
$ws.Range("E58").Value = 523
$ws.Range("H58").Value = 6

# Row 59: Serbia -> Egipto
$ws.Range("A59").Value = "Egipto"
$ws.Range("B59").Value = 536
$ws.Range("D59").Value = 116
$ws.Range("E59").Value = 390
$ws.Range("F59").Value = 0
$ws.Range("H59").Value = 30

# Row 60: Hong Kong -> Serbia
$ws.Range("A60").Value = "Serbia"
$ws.Range("B60").Value = 528
$ws.Range("D60").Value = 42
$ws.Range("E60").Value = 478
$ws.Range("F60").Value = 25
$ws.Range("H60").Value = 8

# Row 136: Guayana Francesa -> Zambia
$ws.Range("A136").Value = "Zambia"
$ws.Range("C136").Value = 6
$ws.Range("D136").Value = 0
$ws.Range("E136").Value = 28

# Row 137: Madagascar -> Guayana Francesa
$ws.Range("A137").Value = "Guayana Francesa"
$ws.Range("B137").Value = 28
$ws.Range("D137").Value = 6
$ws.Range("E137").Value = 22

# Row 138: Barbados -> Madagascar
$ws.Range("A138").Value = "Madagascar"

# Row 139: Togo -> Barbados
$ws.Range("A139").Value = "Barbados"
$ws.Range("B139").Value = 26
$ws.Range("D139").Value = 0
$ws.Range("E139").Value = 26
$ws.Range("H139").Value = 0

# Row 140: Uganda -> Togo
$ws.Range("A140").Value = "Togo"
$ws.Range("B140").Value = 25
$ws.Range("D140").Value = 1
$ws.Range("H140").Value = 1

# Row 141: Zambia -> Uganda
$ws.Range("A141").Value = "Uganda"
$ws.Range("B141").Value = 23
$ws.Range("E141").Value = 23
